$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New admin-facing columns: to_show_late / to_show_long ---
# Order matters for shared-string table layout, so update string-valued
# cells in the same order the target workbook lists them.
$ws.Range("E6").Value = "Bathinda"
$ws.Range("C6").Value = "Raghav"
$ws.Range("I5").Value = "abc"
$ws.Range("I6").Value = "xyz"
$ws.Range("AB1").Value = "to_show_late"
$ws.Range("AC1").Value = "to_show_long"

# New boolean data columns, defaulted to FALSE for every student row
$ws.Range("AB2").Value = $false
$ws.Range("AC2").Value = $false
$ws.Range("AB3").Value = $false
$ws.Range("AC3").Value = $false
$ws.Range("AB4").Value = $false
$ws.Range("AC4").Value = $false
$ws.Range("AB5").Value = $false
$ws.Range("AC5").Value = $false
$ws.Range("AB6").Value = $false
$ws.Range("AC6").Value = $false

# Reset a few existing leave-application flags back to FALSE
$ws.Range("N2").Value = $false
$ws.Range("Z2").Value = $false
$ws.Range("O3").Value = $false

# Update the active sheet view / selection to reflect the admin's working area
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("N10").Select() | Out-Null
